$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows("8:12").Delete()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B6").Value = 50
$ws2.Rows("7:7").Delete()
